$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text values: force Text format first so
# Excel keeps them as literal strings (preserving formats like trailing zeros)
# instead of auto-converting to numbers.
$textGuardCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D17", "D19", "D22", "D23", "D25", "D27", "D29", "D30", "D31", "D34", "D35", "D36", "D37", "D38", "D45", "D47", "D48")
foreach ($addr in $textGuardCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '42.952.65'
$ws.Range("E2").Value = '  -1.39%  '
$ws.Range("D3").Value = '2.341.54'
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '307.01'
$ws.Range("E5").Value = '  -1.41%  '
$ws.Range("D6").Value = '100.77'
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("D7").Value = '0.510'
$ws.Range("E7").Value = '  -5.03%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.512'
$ws.Range("E9").Value = '  -3.95%  '
$ws.Range("D10").Value = '34.94'
$ws.Range("E10").Value = '  -2.67%  '
$ws.Range("D11").Value = '52.23'
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").Value = '0.0801'
$ws.Range("E12").Value = '  -2.09%  '
$ws.Range("E13").Value = '  -0.46%  '
$ws.Range("D14").Value = '6.82'
$ws.Range("E14").Value = '  -3.02%  '
$ws.Range("D15").Value = '15.84'
$ws.Range("E15").Value = '  +5.16%  '
$ws.Range("D16").Value = '2.330.42'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '0.804'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '42.860.78'
$ws.Range("E18").Value = '  -1.37%  '
$ws.Range("D19").Value = '6.24'
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("E21").Value = '  -6.26%  '
$ws.Range("D22").Value = '67.87'
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").Value = '236.98'
$ws.Range("E24").Value = '  -2.06%  '
$ws.Range("D25").Value = '2.57'
$ws.Range("E25").Value = '  -2.11%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("D27").Value = '25.52'
$ws.Range("E27").Value = '  +2.67%  '
$ws.Range("E28").Value = '  +9.40%  '
$ws.Range("D29").Value = '35.03'
$ws.Range("E29").Value = '  -5.06%  '
$ws.Range("D30").Value = '9.36'
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("D31").Value = '159.84'
$ws.Range("E31").Value = '  -4.80%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  -3.20%  '
$ws.Range("D34").Value = '4.67'
$ws.Range("E34").Value = '  +8.00%  '
$ws.Range("D35").Value = '2.46'
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("D36").Value = '0.0728'
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("D37").Value = '17.36'
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("D38").Value = '2.97'
$ws.Range("E38").Value = '  -4.80%  '
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("E40").Value = '  -3.37%  '
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("D43").Value = '2.028.58'
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").Value = '18.83'
$ws.Range("E45").Value = '  -2.97%  '
$ws.Range("E46").Value = '  +3.97%  '
$ws.Range("D47").Value = '2.95'
$ws.Range("E47").Value = '  -1.55%  '
$ws.Range("D48").Value = '56.31'
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("E49").Value = '  -1.29%  '
$ws.Range("D50").Value = '2.566.83'
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("E51").Value = '  +2.33%  '
